# Revert "new changes in ops (ordercreation & orderpage & order form)"
#
# The prior commit had added a new order row (row 3), two new columns
# ("Typist" / "Typist QC" in columns E:F) and widened/renumbered a couple
# of other columns. This script reverts all of that, restoring the sheet
# to a single header row + single data row with the "Typist"/"Typist QC"
# columns removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second data row (order "Be18-002") that was added.
$ws.Rows("3").Delete()

# Remove the "Typist" / "Typist QC" columns that were added (columns E:F);
# this shifts Client/Lob/Process/... back left into their original slots.
$ws.Columns("E:F").Delete()

# Restore the custom (non bestFit) width of the "Emp ID-Order Assigned"
# column and the width of the "Product Name" column.
$ws.Columns("C").ColumnWidth = 35.498697916666664
$ws.Columns("H").ColumnWidth = 28.721354166666668

# Restore the previous selection.
$ws.Range("C10").Select()
